# Update crypto price/volume figures per the Aug 31 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.082.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.525.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.81%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.523.08'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.348'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.943.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.919.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.509.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('E23').Value = '  +2.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.423'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('E26').Value = '  -1.70%  '
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.37%  '
$ws.Range('E29').Value = '  -2.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0774'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '168.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.45%  '
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.12'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('E38').Value = '  -3.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.68'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.813'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '284.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '132.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0508'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0221'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.37'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.72%  '
